$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.282.76'
$ws.Range('E2').Value = '  -0.04%  '

$ws.Range('D3').Value = '1.914.25'
$ws.Range('E3').Value = '  -0.87%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.49%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7406'
$ws.Range('E5').Value = '  -3.46%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.81'
$ws.Range('E6').Value = '  -1.80%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.25%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3150'
$ws.Range('E8').Value = '  -2.26%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.16'
$ws.Range('E9').Value = '  -4.50%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07014'
$ws.Range('E10').Value = '  -1.12%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7825'
$ws.Range('E11').Value = '  -0.94%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07970'
$ws.Range('E12').Value = '  -0.32%  '

$ws.Range('D13').Value = '1.898.63'
$ws.Range('E13').Value = '  -1.77%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.307'
$ws.Range('E14').Value = '  -1.29%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.16'
$ws.Range('E15').Value = '  -2.71%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.38'
$ws.Range('E16').Value = '  -2.49%  '

$ws.Range('D17').Value = '30.229.19'
$ws.Range('E17').Value = '  -0.28%  '

$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.870'
$ws.Range('E18').Value = '  +1.23%  '

$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.71'
$ws.Range('E19').Value = '  -3.13%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007860'
$ws.Range('E20').Value = '  -1.99%  '

$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.11%  '

$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.155.91'
$ws.Range('E22').Value = '  -1.49%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  +0.49%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.679'
$ws.Range('E24').Value = '  -2.15%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.534'
$ws.Range('E25').Value = '  -0.64%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.30'
$ws.Range('E26').Value = '  +0.34%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.04'
$ws.Range('E27').Value = '  -0.38%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1277'
$ws.Range('E28').Value = '  -6.12%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.129'
$ws.Range('E29').Value = '  -8.26%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.358'
$ws.Range('E30').Value = '  -0.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.549'
$ws.Range('E31').Value = '  +1.39%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.336'
$ws.Range('E32').Value = '  -2.32%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.090'
$ws.Range('E33').Value = '  -1.24%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05228'
$ws.Range('E34').Value = '  +1.22%  '

$ws.Range('E35').Value = '  +1.16%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7538'
$ws.Range('E36').Value = '  +0.15%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.759'
$ws.Range('E37').Value = '  -0.38%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01952'
$ws.Range('E38').Value = '  -0.64%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.793'
$ws.Range('E39').Value = '  -0.17%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.420'
$ws.Range('E40').Value = '  +0.01%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '76.12'
$ws.Range('E41').Value = '  -2.95%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4501'
$ws.Range('E42').Value = '  -0.40%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.946'
$ws.Range('E43').Value = '  -2.38%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9993'
$ws.Range('E44').Value = '  -0.02%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.793'
$ws.Range('E45').Value = '  +3.13%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.8322'
$ws.Range('E46').Value = '  -0.46%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.947'
$ws.Range('E47').Value = '  +1.43%  '

$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.24'
$ws.Range('E48').Value = '  -1.41%  '

$ws.Range('D49').Value = '2.077.19'
$ws.Range('E49').Value = '  -0.65%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '37.33'
$ws.Range('E50').Value = '  -0.61%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1219'
$ws.Range('E51').Value = '  +1.95%  '
